$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.853.19'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '1.564.01'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''205.95'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('D8').Value = '''21.78'
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('D10').Value = '''0.0584'
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('D12').Value = '1.786.12'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = '1.563.08'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').Value = '''3.73'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('D16').Value = '26.861.97'
$ws.Range('D17').Value = '''61.27'
$ws.Range('E17').Value = '  -2.58%  '
$ws.Range('D18').Value = '''215.08'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('E20').Value = '  -1.13%  '
$ws.Range('D22').Value = '''4.12'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = '''9.20'
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('D25').Value = '''153.50'
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('E26').Value = '  +2.02%  '
$ws.Range('D27').Value = '''14.91'
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('D30').Value = '''0.0467'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('E31').Value = '  -3.70%  '
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('D33').Value = '1.405.85'
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E36').Value = '  -0.48%  '
$ws.Range('D37').Value = '''0.917'
$ws.Range('E37').Value = '  -2.85%  '
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('D39').Value = '''0.528'
$ws.Range('E39').Value = '  +1.71%  '
$ws.Range('E40').Value = '  -0.77%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  +0.49%  '
$ws.Range('D43').Value = '''5.46'
$ws.Range('E43').Value = '  +4.51%  '
$ws.Range('E44').Value = '  -1.45%  '
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('D46').Value = '''63.36'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('D47').Value = '1.699.55'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').Value = '''0.0505'
$ws.Range('E49').Value = '  +2.50%  '
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('E51').Value = '  +0.51%  '
